# Update market-price derived figures (currentAveragePrice*, Leve cost/profit columns)
# across the Leve tracking sheets, per the latest scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9999.5
$ws.Range("I32").Value = 9999.5
$ws.Range("K32").Value = 9999.5
$ws.Range("M32").Value = -9673.5
$ws.Range("H62").Value = 10211
$ws.Range("I62").Value = 2460.8572
$ws.Range("J62").Value = 16992.375
$ws.Range("K62").Value = 2460.8572
$ws.Range("L62").Value = 16992.375
$ws.Range("M62").Value = -1836.8572
$ws.Range("N62").Value = -18240.375
$ws.Range("H65").Value = 10211
$ws.Range("I65").Value = 2460.8572
$ws.Range("J65").Value = 16992.375
$ws.Range("K65").Value = 12304.286
$ws.Range("L65").Value = 84961.875
$ws.Range("M65").Value = -9184.286
$ws.Range("N65").Value = -91201.875
$ws.Range("H125").Value = 2418.2727
$ws.Range("I125").Value = 1302
$ws.Range("J125").Value = 2836.875
$ws.Range("K125").Value = 11718
$ws.Range("L125").Value = 25531.875
$ws.Range("M125").Value = -9258
$ws.Range("N125").Value = -30451.875
$ws.Range("H135").Value = 1019.381
$ws.Range("I135").Value = 882.6667
$ws.Range("K135").Value = 7944.0003
$ws.Range("M135").Value = -5409.0003
$ws.Range("H138").Value = 2756.6938
$ws.Range("I138").Value = 1227.7368
$ws.Range("K138").Value = 3683.2104
$ws.Range("M138").Value = 1456.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 150000
$ws.Range("J7").Value = 150000
$ws.Range("L7").Value = 150000
$ws.Range("N7").Value = -150228
$ws.Range("H32").Value = 5972464.5
$ws.Range("I32").Value = 7160322.5
$ws.Range("J32").Value = 33172.43
$ws.Range("K32").Value = 7160322.5
$ws.Range("L32").Value = 33172.43
$ws.Range("M32").Value = -7160035.5
$ws.Range("N32").Value = -33746.43
$ws.Range("H45").Value = 33336230
$ws.Range("I45").Value = 35717104
$ws.Range("K45").Value = 35717104
$ws.Range("M45").Value = -35716727
$ws.Range("H61").Value = 50006084
$ws.Range("I61").Value = 45458732
$ws.Range("J61").Value = 62511308
$ws.Range("K61").Value = 45458732
$ws.Range("L61").Value = 62511308
$ws.Range("M61").Value = -45458520
$ws.Range("N61").Value = -62511732
$ws.Range("H74").Value = 9268934
$ws.Range("J74").Value = 16261.077
$ws.Range("L74").Value = 16261.077
$ws.Range("N74").Value = -18009.077
$ws.Range("H77").Value = 9268934
$ws.Range("J77").Value = 16261.077
$ws.Range("L77").Value = 81305.38499999999
$ws.Range("N77").Value = -90041.38499999999
$ws.Range("H88").Value = 3561.818
$ws.Range("I88").Value = 2916.2
$ws.Range("J88").Value = 4099.8335
$ws.Range("K88").Value = 2916.2
$ws.Range("L88").Value = 4099.8335
$ws.Range("M88").Value = -2510.2
$ws.Range("N88").Value = -4911.8335
$ws.Range("H91").Value = 3561.818
$ws.Range("I91").Value = 2916.2
$ws.Range("J91").Value = 4099.8335
$ws.Range("K91").Value = 2916.2
$ws.Range("L91").Value = 4099.8335
$ws.Range("M91").Value = -1512.2
$ws.Range("N91").Value = -6907.8335
$ws.Range("H97").Value = 1461.1364
$ws.Range("I97").Value = 1461.1364
$ws.Range("K97").Value = 1461.1364
$ws.Range("M97").Value = -965.1364000000001
$ws.Range("H102").Value = 3261.3333
$ws.Range("I102").Value = 4027.9524
$ws.Range("K102").Value = 4027.9524
$ws.Range("M102").Value = -2405.9524
$ws.Range("H132").Value = 7096096
$ws.Range("I132").Value = 9261051
$ws.Range("K132").Value = 27783153
$ws.Range("M132").Value = -27780623
$ws.Range("H136").Value = 50006084
$ws.Range("I136").Value = 45458732
$ws.Range("J136").Value = 62511308
$ws.Range("K136").Value = 136376196
$ws.Range("L136").Value = 187533924
$ws.Range("M136").Value = -136373646
$ws.Range("N136").Value = -187539024
$ws.Range("H139").Value = 100715
$ws.Range("J139").Value = 100715
$ws.Range("L139").Value = 100715
$ws.Range("N139").Value = -110995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3332.2812
$ws.Range("I86").Value = 3122.5833
$ws.Range("K86").Value = 3122.5833
$ws.Range("M86").Value = -1999.5833
$ws.Range("H89").Value = 3332.2812
$ws.Range("I89").Value = 3122.5833
$ws.Range("K89").Value = 15612.9165
$ws.Range("M89").Value = -9996.916499999999
$ws.Range("H134").Value = 169075.11
$ws.Range("I134").Value = 1669.7872
$ws.Range("K134").Value = 5009.3616
$ws.Range("M134").Value = -2474.3616

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1180278.2
$ws.Range("I31").Value = 16531.1
$ws.Range("K31").Value = 16531.1
$ws.Range("M31").Value = -16236.1
$ws.Range("H34").Value = 1180278.2
$ws.Range("I34").Value = 16531.1
$ws.Range("K34").Value = 16531.1
$ws.Range("M34").Value = -16329.1
$ws.Range("H86").Value = 5784.4287
$ws.Range("J86").Value = 5631.25
$ws.Range("L86").Value = 5631.25
$ws.Range("N86").Value = -7877.25
$ws.Range("H89").Value = 5784.4287
$ws.Range("J89").Value = 5631.25
$ws.Range("L89").Value = 28156.25
$ws.Range("N89").Value = -39388.25
$ws.Range("H141").Value = 176926.69
$ws.Range("I141").Value = 42195.25
$ws.Range("J141").Value = 206867
$ws.Range("K141").Value = 42195.25
$ws.Range("L141").Value = 206867
$ws.Range("M141").Value = -37015.25
$ws.Range("N141").Value = -217227

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1847.9048
$ws.Range("J5").Value = 1930.25
$ws.Range("L5").Value = 5790.75
$ws.Range("N5").Value = -6014.75
$ws.Range("H34").Value = 1002.8333
$ws.Range("I34").Value = 403.4
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1210.2
$ws.Range("L34").Value = 12000
$ws.Range("M34").Value = -1126.2
$ws.Range("N34").Value = -12168
$ws.Range("H132").Value = 2579.4
$ws.Range("J132").Value = 2762.9092
$ws.Range("L132").Value = 24866.1828
$ws.Range("N132").Value = -29926.1828
$ws.Range("H135").Value = 1847.9048
$ws.Range("J135").Value = 1930.25
$ws.Range("L135").Value = 17372.25
$ws.Range("N135").Value = -22442.25
$ws.Range("H139").Value = 3072.303
$ws.Range("I139").Value = 1874.65
$ws.Range("K139").Value = 5623.950000000001
$ws.Range("M139").Value = -483.9500000000007
$ws.Range("H140").Value = 70454.25999999999
$ws.Range("I140").Value = 105429.35
$ws.Range("K140").Value = 316288.05
$ws.Range("M140").Value = -311108.05

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 52635370
$ws.Range("I132").Value = 62504012
$ws.Range("J132").Value = 2611.6667
$ws.Range("K132").Value = 187512036
$ws.Range("L132").Value = 7835.000100000001
$ws.Range("M132").Value = -187509506
$ws.Range("N132").Value = -12895.0001
$ws.Range("H141").Value = 33333
$ws.Range("J141").Value = 33333
$ws.Range("L141").Value = 33333
$ws.Range("N141").Value = -43693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 589061.2
$ws.Range("I7").Value = 10380
$ws.Range("J7").Value = 1384747.8
$ws.Range("K7").Value = 10380
$ws.Range("L7").Value = 1384747.8
$ws.Range("M7").Value = -10268
$ws.Range("N7").Value = -1384971.8
$ws.Range("H40").Value = 6000
$ws.Range("I40").Value = 6000
$ws.Range("K40").Value = 6000
$ws.Range("M40").Value = -5864
$ws.Range("H46").Value = 4029.0908
$ws.Range("I46").Value = 2360
$ws.Range("J46").Value = 5420
$ws.Range("K46").Value = 2360
$ws.Range("L46").Value = 5420
$ws.Range("M46").Value = -2172
$ws.Range("N46").Value = -5796
$ws.Range("H100").Value = 3618.0908
$ws.Range("I100").Value = 2883.1667
$ws.Range("J100").Value = 4500
$ws.Range("K100").Value = 2883.1667
$ws.Range("L100").Value = 4500
$ws.Range("M100").Value = -2342.1667
$ws.Range("N100").Value = -5582
$ws.Range("H122").Value = 5464.3794
$ws.Range("I122").Value = 4702.1
$ws.Range("J122").Value = 7158.3335
$ws.Range("K122").Value = 14106.3
$ws.Range("L122").Value = 21475.0005
$ws.Range("M122").Value = -11656.3
$ws.Range("N122").Value = -26375.0005
$ws.Range("H126").Value = 589061.2
$ws.Range("I126").Value = 10380
$ws.Range("J126").Value = 1384747.8
$ws.Range("K126").Value = 31140
$ws.Range("L126").Value = 4154243.4
$ws.Range("M126").Value = -28670
$ws.Range("N126").Value = -4159183.4
$ws.Range("H132").Value = 3928.1875
$ws.Range("I132").Value = 5456.5713
$ws.Range("J132").Value = 2739.4443
$ws.Range("K132").Value = 16369.7139
$ws.Range("L132").Value = 8218.332900000001
$ws.Range("M132").Value = -13839.7139
$ws.Range("N132").Value = -13278.3329
$ws.Range("H136").Value = 83554.375
$ws.Range("I136").Value = 21492.875
$ws.Range("K136").Value = 64478.625
$ws.Range("M136").Value = -61928.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1335.6316
$ws.Range("I100").Value = 1683.5714
$ws.Range("J100").Value = 1132.6666
$ws.Range("K100").Value = 3367.1428
$ws.Range("L100").Value = 2265.3332
$ws.Range("M100").Value = -2826.1428
$ws.Range("N100").Value = -3347.3332
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
